# Add new match-day rows to the "Partidos" sheet (rows 496-501, date 2025-10-25 / serial 45955)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Partidos")

$fecha = 45955

$rows = @(
    @{ Row=496; Jugador="Gember Marin Sarria";          Equipo="Amarillo"; Posicion="Arquero";       Goles=0; Autogoles=0; Arquero=$true;  GolesRecibidos=0; Amarillas=0; Rojas=0; Asistencias=0; Penales=0 },
    @{ Row=497; Jugador="Edwin Casas";                   Equipo="Azul";     Posicion="Arquero";       Goles=0; Autogoles=0; Arquero=$true;  GolesRecibidos=3; Amarillas=0; Rojas=0; Asistencias=0; Penales=0 },
    @{ Row=498; Jugador="Alexander Uribe";                Equipo="Amarillo"; Posicion="Mediocampista"; Goles=1; Autogoles=0; Arquero=$false; GolesRecibidos=0; Amarillas=0; Rojas=0; Asistencias=0; Penales=0 },
    @{ Row=499; Jugador="Armando Murillo";                Equipo="Amarillo"; Posicion="Defensa";       Goles=1; Autogoles=0; Arquero=$false; GolesRecibidos=0; Amarillas=0; Rojas=0; Asistencias=0; Penales=0 },
    @{ Row=500; Jugador="Carlos Fernando Valencia";       Equipo="Amarillo"; Posicion="Delantero";     Goles=1; Autogoles=0; Arquero=$false; GolesRecibidos=0; Amarillas=0; Rojas=0; Asistencias=1; Penales=0 },
    @{ Row=501; Jugador="Bryan Andres Burgos";            Equipo="Amarillo"; Posicion="Mediocampista"; Goles=0; Autogoles=0; Arquero=$false; GolesRecibidos=0; Amarillas=0; Rojas=0; Asistencias=1; Penales=0 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $fecha
    $ws.Cells.Item($row, 2).Value = $r.Jugador
    $ws.Cells.Item($row, 3).Value = $r.Equipo
    $ws.Cells.Item($row, 4).Value = $r.Posicion
    $ws.Cells.Item($row, 5).Value = $r.Goles
    $ws.Cells.Item($row, 6).Value = $r.Autogoles
    $ws.Cells.Item($row, 7).Value = $r.Arquero
    $ws.Cells.Item($row, 8).Value = $r.GolesRecibidos
    $ws.Cells.Item($row, 9).Value = $r.Amarillas
    $ws.Cells.Item($row, 10).Value = $r.Rojas
    $ws.Cells.Item($row, 11).Value = $r.Asistencias
    $ws.Cells.Item($row, 12).Value = $r.Penales
}

# Update view state to match the author's final selection/scroll position
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 483
$ws.Range("A502").Select()
